# Automatische test-sync: 2025-08-03 18:46:50
# Adds a new test-mail row (#16) to the "Logs" sheet and updates the
# "Dashboard" category-count table + chart source range accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 44 with the new test-mail entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 44
$logs.Cells.Item($newRow, 1).Value  = "Wil je dit even doorsturen?"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #16: Wil je dit even doorsturen?"
$logs.Cells.Item($newRow, 4).Value  = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-03 18:46:43"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 43 to
# row 44 so the new row participates in the same highlighting rules.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`43")
    $newRange = $logs.Range("$col`2:$col`44")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: "Planning / Afspraak" now ties with "Overig" at
#    11 occurrences each, so the two rows swap places in the table
#    feeding the bar chart.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Planning / Afspraak"
$dashboard.Range("B2").Value = 11
$dashboard.Range("A3").Value = "Overig"
$dashboard.Range("B3").Value = 11
